$wb = $excel.ActiveWorkbook

# Update the selection left behind on the original "UserModel" sheet.
$userModelSheet = $wb.Worksheets.Item("UserModel")
$userModelSheet.Activate()
[void]$userModelSheet.Range("A35").Select()

# Add the new worksheet "Sheet1" after the existing "UserModel" sheet.
$ws = $wb.Worksheets.Add([System.Type]::Missing, $userModelSheet)
$ws.Name = "Sheet1"

# Row 4
$ws.Range("A4").Value = "userName"
$ws.Range("B4").Value = "test"
$ws.Range("C4").Value = "test2"
$ws.Range("D4").Value = "test3"

# Row 5
$ws.Range("A5").Value = "login"
$ws.Range("B5").Value = "test"

# Row 6
$ws.Range("A6").Value = "pass"
$ws.Range("B6").Value = "test"

# Row 7
$ws.Range("A7").Value = "userNameNick"
$ws.Range("B7").Value = "Test"
$ws.Range("C7").Value = "Test2"
$ws.Range("D7").Value = "Test3"

# Row 8
$ws.Range("A8").Value = "Skills"
$ws.Range("B8").Value = "A"
$ws.Range("C8").Value = "A"
$ws.Range("D8").Value = "A"

# Row 9
$ws.Range("A9").Value = "Language"
$ws.Range("B9").Value = "En"
$ws.Range("C9").Value = "En, Ch"
$ws.Range("D9").Value = "En"

# Row 10
$ws.Range("A10").Value = "Country"
$ws.Range("B10").Value = "US"
$ws.Range("C10").Value = "US"
$ws.Range("D10").Value = "US"

# Row 11
$ws.Range("A11").Value = "Gender"
$ws.Range("B11").Value = "male"
$ws.Range("C11").Value = "female"
$ws.Range("D11").Value = "male"

$ws.Columns.Item(1).ColumnWidth = 14

# Set the active sheet/selection/zoom to match the target state
$ws.Activate()
$excel.ActiveWindow.Zoom = 189
[void]$ws.Range("C8").Select()
